$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value updates (prices in column D, 1h volume % in column E,
# plus a couple of coin-name/link swaps in columns B/C).
$ws.Range('D2').Value = '67.510.49'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '3.494.28'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''606.50'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = '''150.72'
$ws.Range('E6').Value = '  +0.96%  '
$ws.Range('D7').Value = '3.492.81'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.485'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('E11').Value = '  +6.35%  '
$ws.Range('D12').Value = '''0.429'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '''32.40'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').Value = '4.094.37'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.499.99'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '67.468.68'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').Value = '''6.52'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').Value = '''15.45'
$ws.Range('E20').Value = '  +2.09%  '
$ws.Range('D21').Value = '''9.73'
$ws.Range('E21').Value = '  +5.93%  '
$ws.Range('D22').Value = '''445.47'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '''0.628'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').Value = '''77.63'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '3.642.27'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '''0.0000126'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').Value = '''8.80'
$ws.Range('E28').Value = '  +6.19%  '
$ws.Range('D29').Value = '''10.02'
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('D31').Value = '''1.64'
$ws.Range('E31').Value = '  +6.87%  '
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '''25.61'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').Value = '''6.14'
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''1.86'
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.491.40'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').Value = '''7.96'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  +5.12%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '''174.77'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').Value = '''0.0891'
$ws.Range('E43').Value = '  +2.44%  '
$ws.Range('D44').Value = '''5.47'
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').Value = '''29.77'
$ws.Range('E46').Value = '  +8.92%  '
$ws.Range('D47').Value = '''46.49'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('E48').Value = '  +3.70%  '
$ws.Range('D49').Value = '''7.62'
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').Value = '''0.252'
$ws.Range('E51').Value = '  +3.36%  '

# The apostrophe-prefixed numeric-looking price cells above pick up an
# auto-applied "quote prefix" style. Re-stamp each with a plain text-style
# cell's style (one cell at a time -- a comma-joined Range union only
# restyles its first area) so the cells stay unstyled, like every other
# data cell in this sheet.
$plainStyle = $ws.Range("D4").Style
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").Style = $plainStyle
$ws.Range("D9").Style = $plainStyle
$ws.Range("D12").Style = $plainStyle
$ws.Range("D13").Style = $plainStyle
$ws.Range("D19").Style = $plainStyle
$ws.Range("D20").Style = $plainStyle
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").Style = $plainStyle
$ws.Range("D23").Style = $plainStyle
$ws.Range("D24").Style = $plainStyle
$ws.Range("D27").Style = $plainStyle
$ws.Range("D28").Style = $plainStyle
$ws.Range("D29").Style = $plainStyle
$ws.Range("D31").Style = $plainStyle
$ws.Range("D34").Style = $plainStyle
$ws.Range("D35").Style = $plainStyle
$ws.Range("D36").Style = $plainStyle
$ws.Range("D38").Style = $plainStyle
$ws.Range("D42").Style = $plainStyle
$ws.Range("D43").Style = $plainStyle
$ws.Range("D44").Style = $plainStyle
$ws.Range("D46").Style = $plainStyle
$ws.Range("D47").Style = $plainStyle
$ws.Range("D49").Style = $plainStyle
$ws.Range("D51").Style = $plainStyle
